$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated B and C values for rows 2-47 (time 0-45)
$data = @(
    @{ Row = 2; B = 2.613752437523362; C = 0.7592657640396985 },
    @{ Row = 3; B = 2.716536614101855; C = 1.460595639285859 },
    @{ Row = 4; B = 7.041282202434871; C = 2.304059489169519 },
    @{ Row = 5; B = 12.03599109472344; C = 3.173738742024204 },
    @{ Row = 6; B = 14.6921863375777; C = 4.001888878368169 },
    @{ Row = 7; B = 21.57270026272947; C = 4.887993696444225 },
    @{ Row = 8; B = 21.68642128141905; C = 5.959203634897216 },
    @{ Row = 9; B = 22.12642842980381; C = 6.876654006919253 },
    @{ Row = 10; B = 22.43733789614787; C = 7.845604822855141 },
    @{ Row = 11; B = 22.47942909118693; C = 8.645824804633929 },
    @{ Row = 12; B = 23.68583780630273; C = 9.651851359695083 },
    @{ Row = 13; B = 24.6012085290772; C = 10.37624620878069 },
    @{ Row = 14; B = 26.95666560595685; C = 11.0979363918407 },
    @{ Row = 15; B = 28.51084660616283; C = 11.94586713835366 },
    @{ Row = 16; B = 28.64542511507138; C = 12.88993087357854 },
    @{ Row = 17; B = 30.00697626199472; C = 13.73978214927143 },
    @{ Row = 18; B = 30.18160745226803; C = 14.9431245371633 },
    @{ Row = 19; B = 32.29221835264032; C = 15.84098931660236 },
    @{ Row = 20; B = 36.68977613810812; C = 16.71748676836768 },
    @{ Row = 21; B = 39.24058746232968; C = 17.65458176062602 },
    @{ Row = 22; B = 42.45769340898283; C = 18.38999985699175 },
    @{ Row = 23; B = 42.72038277164717; C = 19.29597632478046 },
    @{ Row = 24; B = 47.12024772749498; C = 20.06302809841178 },
    @{ Row = 25; B = 47.18983974099818; C = 20.84046645116023 },
    @{ Row = 26; B = 47.45139448397205; C = 21.73050188173235 },
    @{ Row = 27; B = 57.08002863459992; C = 22.53922142054101 },
    @{ Row = 28; B = 57.19176874318924; C = 23.46843318447786 },
    @{ Row = 29; B = 57.75826233793363; C = 24.32965619625855 },
    @{ Row = 30; B = 59.94346593310251; C = 25.51196157872716 },
    @{ Row = 31; B = 62.82205560266711; C = 26.38783686786138 },
    @{ Row = 32; B = 62.8817530065571; C = 27.30969612916372 },
    @{ Row = 33; B = 72.54148468743868; C = 28.040704622941 },
    @{ Row = 34; B = 72.97501644923719; C = 28.84470881218256 },
    @{ Row = 35; B = 73.06838145081315; C = 29.75031308741781 },
    @{ Row = 36; B = 75.27802637213172; C = 30.48428657400807 },
    @{ Row = 37; B = 78.20291173745026; C = 31.46848739832284 },
    @{ Row = 38; B = 78.29546518141368; C = 32.30596951609363 },
    @{ Row = 39; B = 78.60360553774709; C = 33.14985736365251 },
    @{ Row = 40; B = 79.72355933701955; C = 33.99770763239682 },
    @{ Row = 41; B = 79.80473470598035; C = 34.88410236901647 },
    @{ Row = 42; B = 94.34852273764618; C = 35.82076498679206 },
    @{ Row = 43; B = 94.47171610910208; C = 36.7427201892187 },
    @{ Row = 44; B = 94.58400135795851; C = 37.50748567541343 },
    @{ Row = 45; B = 95.01755374451145; C = 38.35212812944956 },
    @{ Row = 46; B = 98.58123591599718; C = 39.21071623763321 },
    @{ Row = 47; B = 98.99359252344701; C = 40.56455204753814 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}

# New row 48: time = 46, plus new B/C values
# Copy formatting from the cell above (A47) so A48 keeps the same style (bold, border, centered)
$ws.Cells.Item(47, 1).Copy($ws.Cells.Item(48, 1))
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 99.23171059013568
$ws.Cells.Item(48, 3).Value = 41.37185266170029
